# ADD Fomulario de contacto y arreglos varios
# Adds a new product row (id=4, sku=A004, "Raton") to the Hoja1 sheet,
# duplicating the data from row 4 but giving the product name (C5) an
# underline style, and moves the active selection to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5: duplicate of row 4's data, with id=4 and sku="A004" ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "A004"
$ws.Range("C5").Value = "Raton"
$ws.Range("D5").Value = "un raton super pro"
$ws.Range("E5").Value = "img/raton"

# "preu" column stores its value as text (matches existing rows, e.g. F4),
# so force text formatting while entering the value, then restore the
# number format so the style itself stays plain/default.
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "59.95"
$ws.Range("F5").NumberFormat = "General"

$ws.Range("G5").Value = 20

# Product name on the new row is underlined.
$ws.Range("C5").Font.Underline = 2

# Move the active selection to C5 (was E10).
$ws.Range("C5").Select()
